$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new row of data (row 44): date serial 46001 (12/10/2025), 625, 22, 603
$ws.Range("A44").Value = 46001
$ws.Range("B44").Value = 625
$ws.Range("C44").Value = 22
$ws.Range("D44").Value = 603

# Update the active selection to A44:D44 (matches the sheetView selection change)
$ws.Range("A44:D44").Select()
